$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 corresponds to the Notion "Tháng 8" (August) report row.
# Update the last_edited_time / last_edited_by.id tracking fields
# and fill in the report figures (Dư nợ phát sinh, Chi tiêu, Lũy kế,
# Tổng doanh thu, Đã thanh toán, Số lượng đơn, Thu nợ, Đơn giá).
#
# The source export shares a single string-table entry for
# "2024-07-31T18:24:00.000Z" across several rows' last_edited_time
# column (D3, D4, D5, D7, D13 as well as D10). The upstream edit
# rewrote that shared entry's text in place, so every row pointing at
# it now reads "2024-08-03T03:17:00.000Z" - update all of them to the
# new timestamp so they keep sharing the same value.
$ws.Range("D3").Value = "2024-08-03T03:17:00.000Z"
$ws.Range("D4").Value = "2024-08-03T03:17:00.000Z"
$ws.Range("D5").Value = "2024-08-03T03:17:00.000Z"
$ws.Range("D7").Value = "2024-08-03T03:17:00.000Z"
$ws.Range("D10").Value = "2024-08-03T03:17:00.000Z"
$ws.Range("D13").Value = "2024-08-03T03:17:00.000Z"

$ws.Range("N10").Value = "41cabcaf-915d-46a5-8eff-38727be27269"

$ws.Range("T10").Value = 1000000
$ws.Range("W10").Value = 1857000
$ws.Range("AA10").Value = 9143000
$ws.Range("AE10").Value = 11000000
$ws.Range("AH10").Value = 11000000
$ws.Range("AK10").Value = 2
$ws.Range("AN10").Value = 0
$ws.Range("AQ10").Value = 12000000
